$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 93 - this shifts existing rows 93..201 down to 94..202
$ws.Rows.Item(93).Insert()

# Populate the new row 93 with its data (same constant columns as the rest of the
# dataset, plus the row-specific values taken from the edit)
$ws.Cells.Item(93, 1).Value = 3
$ws.Cells.Item(93, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(93, 3).Value = "Coquimbo"
$ws.Cells.Item(93, 4).Value = 44897
$ws.Cells.Item(93, 5).Value = 5
$ws.Cells.Item(93, 6).Value = 100112030
$ws.Cells.Item(93, 7).Value = "Poroto granado"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 73
$ws.Cells.Item(93, 11).Value = 39000
$ws.Cells.Item(93, 12).Value = 40000
$ws.Cells.Item(93, 13).Value = 39521
$ws.Cells.Item(93, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(93, 16).Value = 1581
$ws.Cells.Item(93, 17).Value = 25
$ws.Cells.Item(93, 18).Value = "Hortaliza"
